$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 11: the "Barker at 1MHz offset x4 ..." pulse now routes its marker
# through marker-2 (column N) instead of marker-1 (column M).
# ---------------------------------------------------------------------------
$ws.Cells.Item(11, 13).Value = 0
$ws.Cells.Item(11, 14).Value = 1
$ws.Cells.Item(11, 16).Value = "Barker at 1MHz offset x4 with external marker 2"

# ---------------------------------------------------------------------------
# Insert a new row at 20 (old row 20 -> row 21). This is the split of the
# former "Barker at 1MHz offset 10dBm 45deg offset x4" pulse (old row 19)
# into a x3 burst (row 19) plus a separate x1-with-marker burst (new row 20).
# ---------------------------------------------------------------------------
$ws.Rows.Item(20).Insert()

# Row 19 now only plays the Barker burst 3 times.
$ws.Cells.Item(19, 12).Value = 3
$ws.Range("B19").Formula = "=B17+K17*(L17+1)"
$ws.Cells.Item(19, 16).Value = "Barker at 1MHz offset 10dBm 45deg offset x3"

# New row 20: single extra Barker pulse with Marker 1 asserted.
$ws.Cells.Item(20, 1).Value = "pdw"
$ws.Range("B20").Formula = "=B18+K18*(L18+1)"
$ws.Cells.Item(20, 3).Value = "Barker"
$ws.Cells.Item(20, 4).Value = 0.000025
$ws.Cells.Item(20, 5).Value = 1000000
$ws.Cells.Item(20, 6).Value = 10
$ws.Cells.Item(20, 7).Value = 45
$ws.Cells.Item(20, 8).Value = 0.00000357
$ws.Cells.Item(20, 9).Value = 6
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0.0001
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 1
$ws.Cells.Item(20, 14).Value = 0
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(20, 16).Value = "Barker at 1MHz offset 10dBm 45deg offset x1 with Marker 1"

# Row 21 (the former row 20 holding the tcdw/EOF marker) keeps referencing
# the row immediately above it, which is now the new row 20.
$ws.Range("B21").Formula = "=B20+K20*(L20+1)"

# ---------------------------------------------------------------------------
# Drop the stray ":w" comment row that used to live at row 27.
# ---------------------------------------------------------------------------
$ws.Rows.Item(28).Delete()

# ---------------------------------------------------------------------------
# View: freeze the header row and scroll the sheet down, with N31 selected
# in the scrollable pane (matches the reviewer's last on-screen position).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("N31").Select()
